# Generate Report for Handback
# Refresh the handoff/handback timestamps for the most recently regenerated
# source file (row 2) on each per-locale status sheet.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-04 18:53:24"
$wsZhCn.Range("K2").Value = "2016-09-04 18:53:40"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-04 18:53:28"
$wsDeDe.Range("K2").Value = "2016-09-04 18:53:47"
